$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.808914661407471
$ws.Range("B1").Value = 1.93298864364624
$ws.Range("C1").Value = 2.087076187133789
$ws.Range("D1").Value = 2.989427328109741
$ws.Range("E1").Value = 3.130656003952026
